$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = 'Volume 30   Number  19'
$ws.Range("C9").Value = 'Report Covering the Week  5/8/2023  Through  5/14/2023'

# --- Data table updates (rows 14-30) ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = '0'
$ws.Range("E14").Value = '***.*'
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -40
$ws.Range("N14").Value = -88.461538461538
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 11.111111111111
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -56.521739130434
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 17.647058823529
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 81
$ws.Range("K16").Value = -13.580246913580
$ws.Range("L16").Value = 40
$ws.Range("M16").Value = -14.634146341463
$ws.Range("N16").Value = -91.150442477876
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -16.216216216216
$ws.Range("I17").Value = 145
$ws.Range("J17").Value = 121
$ws.Range("K17").Value = 19.834710743801
$ws.Range("L17").Value = 55.913978494623
$ws.Range("M17").Value = 40.776699029126
$ws.Range("N17").Value = -63.659147869674
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -4.761904761904
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = -21.917808219178
$ws.Range("L18").Value = 62.857142857142
$ws.Range("M18").Value = -43
$ws.Range("N18").Value = -90.436241610738
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 13.793103448275
$ws.Range("I19").Value = 156
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 4
$ws.Range("L19").Value = 73.333333333333
$ws.Range("M19").Value = 59.183673469387
$ws.Range("N19").Value = -25.358851674641
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = 38.095238095238
$ws.Range("L20").Value = 114.814814814815
$ws.Range("M20").Value = 31.818181818181
$ws.Range("N20").Value = -81.229773462783
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -34.210526315789
$ws.Range("F21").Value = 123
$ws.Range("G21").Value = 115
$ws.Range("H21").Value = 6.956521739130
$ws.Range("I21").Value = 499
$ws.Range("J21").Value = 479
$ws.Range("K21").Value = 4.175365344467
$ws.Range("L21").Value = 63.071895424836
$ws.Range("M21").Value = 13.926940639269
$ws.Range("N21").Value = -78.793030174245
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -83.333333333333
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -46.666666666666
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 52
$ws.Range("J23").Value = 42
$ws.Range("K23").Value = 23.809523809523
$ws.Range("L23").Value = 108
$ws.Range("M23").Value = 126.086956521739
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -37.5
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = -18.055555555555
$ws.Range("I24").Value = 271
$ws.Range("J24").Value = 313
$ws.Range("K24").Value = -13.418530351437
$ws.Range("L24").Value = 11.065573770491
$ws.Range("M24").Value = -7.823129251700
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -4.081632653061
$ws.Range("I25").Value = 217
$ws.Range("J25").Value = 190
$ws.Range("K25").Value = 14.210526315789
$ws.Range("L25").Value = 79.338842975206
$ws.Range("M25").Value = -25.172413793103
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 11
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -38.888888888888
$ws.Range("L26").Value = -8.333333333333
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = -9.090909090909
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = '0'
$ws.Range("E28").Value = '***.*'
$ws.Range("I28").Value = 7
$ws.Range("K28").Value = 16.666666666666
$ws.Range("L28").Value = -53.333333333333
$ws.Range("M28").Value = -69.565217391304
$ws.Range("N28").Value = -92.391304347826
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = '0'
$ws.Range("E29").Value = '***.*'
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = 40
$ws.Range("L29").Value = -53.333333333333
$ws.Range("M29").Value = -61.111111111111
$ws.Range("N29").Value = -91.954022988505
